$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Luca's suggestions: quantities / rows in the BOM table ---

# "cond C" placeholder row (row 24) is no longer needed -> clear it entirely
$ws.Range("C24:E24").ClearContents()
$ws.Range("H24").ClearContents()

# res 10k : quantity 10 -> 6
$ws.Range("C26").Value = 6

# res 145 (50mA) : quantity 4 -> 7
$ws.Range("C30").Value = 7

# "res 330 (50mA)" row (row 31) is no longer needed -> clear it entirely
$ws.Range("C31:F31").ClearContents()
$ws.Range("H31").ClearContents()

# --- placement sur PCB / view state ---
$ws.Activate()
$ws.Range("K24").Select()
